$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New headers for columns F, G, H (Outliers / MAD variants)
$ws.Range("F1").Value = "KNN_Outliers_MAD"
$ws.Range("G1").Value = "SVM_Outliers_MAD"
$ws.Range("H1").Value = "RF_Outliers_MAD"

# Match the header style used by A1:E1 (bold, centered, bordered) by
# copying the format from the existing E1 header cell.
$ws.Range("E1").Copy()
$ws.Range("F1:H1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Fill F2:H21 with boolean FALSE values (matching the diff's <c t="b"><v>0</v></c>)
for ($r = 2; $r -le 21; $r++) {
    $ws.Cells.Item($r, 6).Value = $false
    $ws.Cells.Item($r, 7).Value = $false
    $ws.Cells.Item($r, 8).Value = $false
}
